# Regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 3.182878228561681
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 0.1529057820181812
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("G2").Value = 5.488907176552729

# Row 3
$ws.Range("B3").Value = 0.7287194209349384
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 0.7127328510149897
$ws.Range("E3").Value = 246.9852506941017
$ws.Range("G3").Value = 250.0799394249504

# Row 4
$ws.Range("B4").Value = 1.505614041169197
$ws.Range("C4").Value = 1.65323645889881
$ws.Range("D4").Value = 0.7127328510149897
$ws.Range("E4").Value = 0.4998867070740569
$ws.Range("G4").Value = 4.371470058157054

# Row 5
$ws.Range("B5").Value = 0.1554434735375247
$ws.Range("C5").Value = 0.05231270169004087
$ws.Range("D5").Value = 0.7127328510149897
$ws.Range("E5").Value = 0.4998867070740569
$ws.Range("G5").Value = 1.420375733316612
